# Applies the "Voeg kennis van KA toe aan M23." edit plus the OpenVAS
# paragraph removal described by the diff.

$p = $ppt.ActivePresentation

# --- Slide 13: "M16" tools slide -------------------------------------
# Remove the paragraph about checking the configuration for known
# vulnerabilities (now covered implicitly by the remaining bullet about
# externally used software versions).
$slide13 = $p.Slides.Item(13)
$tb13 = $slide13.Shapes.Item(2)
$tr13 = $tb13.TextFrame.TextRange

$target = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,"
for ($i = 1; $i -le $tr13.Paragraphs().Count; $i++) {
    $para = $tr13.Paragraphs($i, 1)
    # PowerPoint's TextRange.Text includes the trailing paragraph mark
    # (`\r`), so trim it before comparing.
    $paraText = $para.Text.TrimEnd("`r", "`n")
    if ($paraText -eq $target) {
        $para.Delete()
        break
    }
}

# --- Slide 19: "M23" measure slide ------------------------------------
$slide19 = $p.Slides.Item(19)

# Title: add "kennis van en" before "ervaring met de Kwaliteitsaanpak".
$titleShape = $slide19.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleLen = $titleRange.Text.Length
$titleChars = $titleRange.Characters(1, $titleLen)
$titleChars.Text = "M23: Het project zorgt voor de aanwezigheid van kennis van en ervaring met de Kwaliteitsaanpak"

# Body: append a sentence about explaining the Kwaliteitsaanpak to new
# project members who are not yet familiar with it.
$bodyShape = $slide19.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyLen = $bodyRange.Text.Length
$bodyChars = $bodyRange.Characters(1, $bodyLen)
$bodyChars.Text = "De software delivery manager zorgt ervoor dat bij nieuwe projecten wordt gestart met ten minste twee projectleden die bekend zijn met de Kwaliteitsaanpak. Projectleden die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak."
